# Applies the numeric corrections from the scheduled-runner commit.
# The workbook stores plain static values in H:N (no formulas), so each
# affected cell is written directly via Range.Value; a couple of rows also
# gain or lose a trailing LeveProfit cell, handled with ClearContents / Value.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 2726.875
$ws.Range("I6").Value = 802.5
$ws.Range("K6").Value = 2407.5
$ws.Range("M6").Value = -2295.5
$ws.Range("H8").Value = 272.66666
$ws.Range("I8").Value = 272.66666
$ws.Range("K8").Value = 817.9999799999999
$ws.Range("M8").Value = -678.9999799999999
$ws.Range("H31").Value = 54.875
$ws.Range("I31").Value = 54.875
$ws.Range("K31").Value = 164.625
$ws.Range("M31").Value = 65.375
$ws.Range("H62").Value = 8001
$ws.Range("I62").Value = 7599
$ws.Range("K62").Value = 7599
$ws.Range("M62").Value = -6975
$ws.Range("H65").Value = 8001
$ws.Range("I65").Value = 7599
$ws.Range("K65").Value = 37995
$ws.Range("M65").Value = -34875
$ws.Range("H74").Value = 3792.0833
$ws.Range("I74").Value = 3792.0833
$ws.Range("K74").Value = 3792.0833
$ws.Range("M74").Value = -2856.0833
$ws.Range("H77").Value = 3792.0833
$ws.Range("I77").Value = 3792.0833
$ws.Range("K77").Value = 18960.4165
$ws.Range("M77").Value = -14280.4165
$ws.Range("H88").Value = 3100
$ws.Range("J88").Value = 3437.5
$ws.Range("L88").Value = 3437.5
$ws.Range("N88").Value = -4249.5
$ws.Range("H91").Value = 3100
$ws.Range("J91").Value = 3437.5
$ws.Range("L91").Value = 3437.5
$ws.Range("N91").Value = -6245.5
$ws.Range("H98").Value = 2156.0625
$ws.Range("I98").Value = 1499.9333
$ws.Range("J98").Value = 11998
$ws.Range("K98").Value = 1499.9333
$ws.Range("L98").Value = 11998
$ws.Range("M98").Value = -1.933299999999917
$ws.Range("N98").Value = -14994
$ws.Range("H122").Value = 2156.0625
$ws.Range("I122").Value = 1499.9333
$ws.Range("J122").Value = 11998
$ws.Range("K122").Value = 4499.7999
$ws.Range("L122").Value = 35994
$ws.Range("M122").Value = -2049.7999
$ws.Range("N122").Value = -40894
$ws.Range("H132").Value = 2408.1428
$ws.Range("I132").Value = 1977.7556
$ws.Range("J132").Value = 7250
$ws.Range("K132").Value = 5933.266799999999
$ws.Range("L132").Value = 21750
$ws.Range("M132").Value = -3403.266799999999
$ws.Range("N132").Value = -26810
$ws.Range("H135").Value = 2115.3794
$ws.Range("I135").Value = 2082.5
$ws.Range("K135").Value = 18742.5
$ws.Range("M135").Value = -16207.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9159.233
$ws.Range("I32").Value = 3251.3044
$ws.Range("J32").Value = 28571
$ws.Range("K32").Value = 3251.3044
$ws.Range("L32").Value = 28571
$ws.Range("M32").Value = -2964.3044
$ws.Range("N32").Value = -29145
$ws.Range("H63").Value = 6788.9443
$ws.Range("J63").Value = 7939.3076
$ws.Range("L63").Value = 7939.3076
$ws.Range("N63").Value = -9311.3076
$ws.Range("H66").Value = 6788.9443
$ws.Range("J66").Value = 7939.3076
$ws.Range("L66").Value = 39696.538
$ws.Range("N66").Value = -46560.538
$ws.Range("H114").Value = 72999.664
$ws.Range("J114").Value = 72999.664
$ws.Range("L114").Value = 72999.664
$ws.Range("N114").Value = -81677.664
$ws.Range("H132").Value = 2300.7273
$ws.Range("I132").Value = 2116.1042
$ws.Range("K132").Value = 6348.312600000001
$ws.Range("M132").Value = -3818.312600000001
$ws.Range("H133").Value = 100000
$ws.Range("J133").Value = 100000
$ws.Range("L133").Value = 100000
$ws.Range("N133").Value = -105060

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 26317658
$ws.Range("I94").Value = 41668276
$ws.Range("K94").Value = 41668276
$ws.Range("M94").Value = -41667825
$ws.Range("H100").Value = 14955.75
$ws.Range("J100").Value = 14955.75
$ws.Range("L100").Value = 14955.75
$ws.Range("N100").Value = -17119.75

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H43").Value = 111233
$ws.Range("J43").Value = 111233
$ws.Range("L43").Value = 111233
$ws.Range("N43").Value = -111601
$ws.Range("H58").Value = 2072.7273
$ws.Range("I58").Value = 1798.6
$ws.Range("J58").Value = 4814
$ws.Range("K58").Value = 1798.6
$ws.Range("L58").Value = 4814
$ws.Range("M58").Value = -1595.6
$ws.Range("N58").Value = -5220
$ws.Range("H92").Value = 42449.5
$ws.Range("J92").Value = 42449.5
$ws.Range("L92").Value = 42449.5
$ws.Range("N92").Value = -47441.5
$ws.Range("H93").Value = 30000
$ws.Range("I93").Value = 30000
$ws.Range("K93").Value = 30000
$ws.Range("M93").Value = -28128
$ws.Range("H101").Value = 111233
$ws.Range("J101").Value = 111233
$ws.Range("L101").Value = 111233
$ws.Range("N101").Value = -117723
$ws.Range("H102").Value = 58493.668
$ws.Range("J102").Value = 58493.668
$ws.Range("L102").Value = 58493.668
$ws.Range("N102").Value = -63361.668
$ws.Range("H104").Value = 69523
$ws.Range("J104").Value = 69523
$ws.Range("L104").Value = 69523
$ws.Range("N104").Value = -74765
$ws.Range("H105").Value = 3311.6875
$ws.Range("I105").Value = 3108.9092
$ws.Range("J105").Value = 3757.8
$ws.Range("K105").Value = 3108.9092
$ws.Range("L105").Value = 3757.8
$ws.Range("M105").Value = -1361.9092
$ws.Range("N105").Value = -7251.8
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("H109").Value = 68467.2
$ws.Range("J109").Value = 68467.2
$ws.Range("L109").Value = 68467.2
$ws.Range("N109").Value = -70547.2
$ws.Range("H122").Value = 108151.266
$ws.Range("I122").Value = 127314.52
$ws.Range("K122").Value = 381943.56
$ws.Range("M122").Value = -379493.56
$ws.Range("H132").Value = 1991
$ws.Range("I132").Value = 1690.75
$ws.Range("K132").Value = 5072.25
$ws.Range("M132").Value = -2542.25
$ws.Range("H134").Value = 921.8461
$ws.Range("I134").Value = 784.5714
$ws.Range("J134").Value = 1498.4
$ws.Range("K134").Value = 2353.7142
$ws.Range("L134").Value = 4495.200000000001
$ws.Range("M134").Value = 181.2857999999997
$ws.Range("N134").Value = -9565.200000000001
$ws.Range("H136").Value = 2072.7273
$ws.Range("I136").Value = 1798.6
$ws.Range("J136").Value = 4814
$ws.Range("K136").Value = 5395.799999999999
$ws.Range("L136").Value = 14442
$ws.Range("M136").Value = -2845.799999999999
$ws.Range("N136").Value = -19542
$ws.Range("N106").ClearContents()

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 12872.571
$ws.Range("J7").Value = 6696
$ws.Range("L7").Value = 20088
$ws.Range("N7").Value = -20312
$ws.Range("H109").Value = 5590.3335
$ws.Range("I109").Value = 5203.3335
$ws.Range("K109").Value = 15610.0005
$ws.Range("M109").Value = -14570.0005

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H59").Value = 65000
$ws.Range("I59").Value = 30000
$ws.Range("J59").Value = 100000
$ws.Range("K59").Value = 30000
$ws.Range("L59").Value = 100000
$ws.Range("M59").Value = -29417
$ws.Range("N59").Value = -101166
$ws.Range("H80").Value = 5686.533
$ws.Range("I80").Value = 4767.1665
$ws.Range("K80").Value = 4767.1665
$ws.Range("M80").Value = -3769.1665
$ws.Range("H83").Value = 5686.533
$ws.Range("I83").Value = 4767.1665
$ws.Range("K83").Value = 23835.8325
$ws.Range("M83").Value = -18843.8325
$ws.Range("H107").Value = 1353.5454
$ws.Range("I107").Value = 1365.5555
$ws.Range("K107").Value = 1365.5555
$ws.Range("M107").Value = 554.4445000000001
$ws.Range("H113").Value = 6115.0835
$ws.Range("I113").Value = 6398.273
$ws.Range("K113").Value = 6398.273
$ws.Range("M113").Value = -4228.273
$ws.Range("H132").Value = 4881.552
$ws.Range("I132").Value = 4758.204
$ws.Range("K132").Value = 14274.612
$ws.Range("M132").Value = -11744.612

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 9142
$ws.Range("I40").Value = 9142
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 9142
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -9006
$ws.Range("H122").Value = 4995
$ws.Range("I122").Value = 4994
$ws.Range("K122").Value = 14982
$ws.Range("M122").Value = -12532
$ws.Range("N40").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 211.64285
$ws.Range("I107").Value = 213.66667
$ws.Range("K107").Value = 641.00001
$ws.Range("M107").Value = 1278.99999
